$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.025914516299962997
$ws.Range("C2").Value = 0.009059031493961811
$ws.Range("D2").Value = 0.006718204822391272
$ws.Range("E2").Value = 0.005600896663963795
$ws.Range("F2").Value = 0.0000000647206945814105
$ws.Range("G2").Value = 0.0017580282874405384
$ws.Range("J2").Value = 0.12744049727916718
$ws.Range("K2").Value = 1.4362910985946655
